# Update "想去人数" (want-to-go count) values in column F for the
# 展览 (sheet1) and 全部类型 (sheet4) sheets, plus the single row on
# 本地生活 (sheet3), matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 133
$ws1.Range("F3").Value = 127
$ws1.Range("F4").Value = 1271
$ws1.Range("F5").Value = 64
$ws1.Range("F7").Value = 977
$ws1.Range("F8").Value = 939
$ws1.Range("F12").Value = 655
$ws1.Range("F13").Value = 925
$ws1.Range("F14").Value = 1804
$ws1.Range("F15").Value = 3952
$ws1.Range("F16").Value = 1165
$ws1.Range("F17").Value = 112
$ws1.Range("F18").Value = 2606
$ws1.Range("F20").Value = 1080
$ws1.Range("F21").Value = 3584
$ws1.Range("F22").Value = 756
$ws1.Range("F25").Value = 2225
$ws1.Range("F26").Value = 112
$ws1.Range("F27").Value = 841
$ws1.Range("F29").Value = 512
$ws1.Range("F30").Value = 206
$ws1.Range("F31").Value = 110
$ws1.Range("F32").Value = 1339
$ws1.Range("F33").Value = 1953
$ws1.Range("F35").Value = 46
$ws1.Range("F37").Value = 588
$ws1.Range("F38").Value = 279
$ws1.Range("F39").Value = 52
$ws1.Range("F42").Value = 77

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 435

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 435
$ws4.Range("F3").Value = 133
$ws4.Range("F4").Value = 1271
$ws4.Range("F5").Value = 64
$ws4.Range("F6").Value = 977
$ws4.Range("F7").Value = 939
$ws4.Range("F9").Value = 0
$ws4.Range("F14").Value = 925
$ws4.Range("F15").Value = 1804
$ws4.Range("F16").Value = 3952
$ws4.Range("F17").Value = 1165
$ws4.Range("F18").Value = 112
$ws4.Range("F20").Value = 2606
$ws4.Range("F21").Value = 1080
$ws4.Range("F22").Value = 3584
$ws4.Range("F23").Value = 756
$ws4.Range("F27").Value = 2226
$ws4.Range("F31").Value = 112
$ws4.Range("F33").Value = 841
$ws4.Range("F35").Value = 512
$ws4.Range("F36").Value = 206
$ws4.Range("F38").Value = 1339
$ws4.Range("F39").Value = 1953
$ws4.Range("F43").Value = 46
$ws4.Range("F44").Value = 588
$ws4.Range("F45").Value = 279
$ws4.Range("F46").Value = 52
$ws4.Range("F49").Value = 77
